$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Volume/Number text and report-week date range text (rich-text shared strings;
# all runs share identical formatting so merging into a single run is visually
# equivalent to the original multi-run string).
$ws.Range("A8").Value = "Volume 31   Number  10"
$ws.Range("C9").Value = "Report Covering the Week  3/4/2024  Through  3/10/2024"

# Weekly crime-statistics table updates (rows 14-33)
$ws.Range("M14").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("M14").Value = -100
$ws.Range("F15").NumberFormat = "General"
$ws.Range("F15").Value = "'0"
$ws.Range("H15").Value = -100
$ws.Range("L15").Value = -50
$ws.Range("C16").NumberFormat = "General"
$ws.Range("C16").Value = "'0"
$ws.Range("E16").Value = -100
$ws.Range("F16").Value = 9
$ws.Range("G16").Value = 12
$ws.Range("H16").Value = -25
$ws.Range("J16").Value = 28
$ws.Range("K16").Value = -10.714285714285
$ws.Range("L16").Value = -56.896551724137
$ws.Range("M16").Value = 0
$ws.Range("N16").Value = -84.076433121019
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = -50
$ws.Range("F17").Value = 7
$ws.Range("G17").Value = 16
$ws.Range("H17").Value = -56.25
$ws.Range("I17").Value = 35
$ws.Range("J17").Value = 38
$ws.Range("K17").Value = -7.894736842105
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = 40
$ws.Range("N17").Value = -64.285714285714
$ws.Range("C18").NumberFormat = "#,##0"
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 6
$ws.Range("E18").Value = -33.333333333333
$ws.Range("F18").Value = 12
$ws.Range("G18").Value = 20
$ws.Range("H18").Value = -40
$ws.Range("I18").Value = 32
$ws.Range("J18").Value = 53
$ws.Range("K18").Value = -39.622641509434
$ws.Range("L18").Value = -51.515151515151
$ws.Range("M18").Value = -45.762711864406
$ws.Range("N18").Value = -81.818181818181
$ws.Range("C19").Value = 23
$ws.Range("D19").Value = 16
$ws.Range("E19").Value = 43.75
$ws.Range("F19").Value = 57
$ws.Range("G19").Value = 74
$ws.Range("H19").Value = -22.972972972973
$ws.Range("I19").Value = 147
$ws.Range("J19").Value = 170
$ws.Range("K19").Value = -13.529411764705
$ws.Range("L19").Value = -13.529411764705
$ws.Range("M19").Value = 0
$ws.Range("N19").Value = -46.350364963503
$ws.Range("C20").NumberFormat = "#,##0"
$ws.Range("C20").Value = 1
$ws.Range("I20").Value = 4
$ws.Range("K20").Value = -20
$ws.Range("L20").Value = -55.555555555555
$ws.Range("M20").Value = -42.857142857142
$ws.Range("N20").Value = -96
$ws.Range("C21").Value = 30
$ws.Range("D21").Value = 30
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 87
$ws.Range("G21").Value = 126
$ws.Range("H21").Value = -30.952380952381
$ws.Range("I21").Value = 245
$ws.Range("J21").Value = 298
$ws.Range("K21").Value = -17.785234899328
$ws.Range("L21").Value = -28.571428571428
$ws.Range("M21").Value = -7.54716981132
$ws.Range("N21").Value = -70.01223990208
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D22").Value = "'0"
$ws.Range("E22").NumberFormat = "General"
$ws.Range("E22").Value = "'***.*"
$ws.Range("L22").Value = 0
$ws.Range("C23").NumberFormat = "General"
$ws.Range("C23").Value = "'0"
$ws.Range("D23").Value = 3
$ws.Range("E23").Value = -100
$ws.Range("F23").Value = 4
$ws.Range("G23").Value = 10
$ws.Range("H23").Value = -60
$ws.Range("J23").Value = 23
$ws.Range("K23").Value = -52.173913043478
$ws.Range("L23").Value = -64.516129032258
$ws.Range("M23").Value = -52.173913043478
$ws.Range("C24").Value = 32
$ws.Range("D24").Value = 32
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 98
$ws.Range("G24").Value = 119
$ws.Range("H24").Value = -17.647058823529
$ws.Range("I24").Value = 240
$ws.Range("J24").Value = 273
$ws.Range("K24").Value = -12.087912087912
$ws.Range("L24").Value = -35.828877005347
$ws.Range("M24").Value = -18.367346938775
$ws.Range("C25").Value = 18
$ws.Range("D25").Value = 21
$ws.Range("E25").Value = -14.285714285714
$ws.Range("F25").Value = 48
$ws.Range("G25").Value = 75
$ws.Range("H25").Value = -36
$ws.Range("I25").Value = 123
$ws.Range("J25").Value = 171
$ws.Range("K25").Value = -28.070175438596
$ws.Range("L25").Value = -50.403225806451
$ws.Range("C26").Value = 8
$ws.Range("D26").Value = 11
$ws.Range("E26").Value = -27.272727272727
$ws.Range("G26").Value = 29
$ws.Range("H26").Value = -24.137931034482
$ws.Range("I26").Value = 69
$ws.Range("J26").Value = 78
$ws.Range("K26").Value = -11.538461538461
$ws.Range("L26").Value = -15.853658536585
$ws.Range("M26").Value = -16.867469879518
$ws.Range("D27").NumberFormat = "#,##0"
$ws.Range("D27").Value = 1
$ws.Range("E27").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E27").Value = -100
$ws.Range("F27").NumberFormat = "General"
$ws.Range("F27").Value = "'0"
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = -100
$ws.Range("J27").Value = 9
$ws.Range("K27").Value = -66.666666666666
$ws.Range("L27").Value = -62.5
$ws.Range("C28").NumberFormat = "General"
$ws.Range("C28").Value = "'0"
$ws.Range("D28").Value = 3
$ws.Range("E28").Value = -100
$ws.Range("F28").Value = 5
$ws.Range("G28").Value = 6
$ws.Range("H28").Value = -16.666666666666
$ws.Range("J28").Value = 10
$ws.Range("K28").Value = 20
$ws.Range("L28").Value = -36.842105263157
$ws.Range("D29").NumberFormat = "General"
$ws.Range("D29").Value = "'0"
$ws.Range("E29").NumberFormat = "General"
$ws.Range("E29").Value = "'***.*"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("D30").Value = "'0"
$ws.Range("E30").NumberFormat = "General"
$ws.Range("E30").Value = "'***.*"
$ws.Range("I33").Value = 2
